$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: price values in column D are stored as plain text in the source sheet
# (e.g. "42.474.38", "1.00", "0.0000104"). Assigning such strings straight to
# .Value would let Excel's type-inference reinterpret them as numbers, which
# would both lose trailing zeros and mangle values like "42.474.38". Prefixing
# with a leading apostrophe forces Excel to keep them as literal text, exactly
# like typing them into a cell by hand.

# Row 2
$ws.Range("D2").Value = '42.474.38'
$ws.Range("E2").Value = '  +0.98%  '

# Row 3
$ws.Range("D3").Value = '2.285.34'
$ws.Range("E3").Value = '  +0.21%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").Value = '''311.12'
$ws.Range("E5").Value = '  -2.43%  '

# Row 6
$ws.Range("D6").Value = '''103.51'
$ws.Range("E6").Value = '  +2.82%  '

# Row 7
$ws.Range("D7").Value = '''0.622'
$ws.Range("E7").Value = '  -0.67%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").Value = '''0.599'
$ws.Range("E9").Value = '  -0.23%  '

# Row 10
$ws.Range("D10").Value = '''38.90'
$ws.Range("E10").Value = '  +0.50%  '

# Row 11
$ws.Range("D11").Value = '''0.0901'
$ws.Range("E11").Value = '  +0.35%  '

# Row 12
$ws.Range("D12").Value = '''8.24'
$ws.Range("E12").Value = '  +0.52%  '

# Row 13
$ws.Range("E13").Value = '  +1.54%  '

# Row 14
$ws.Range("D14").Value = '''0.979'
$ws.Range("E14").Value = '  +3.09%  '

# Row 15
$ws.Range("D15").Value = '''15.06'
$ws.Range("E15").Value = '  +0.08%  '

# Row 16
$ws.Range("D16").Value = '2.634.46'
$ws.Range("E16").Value = '  +0.25%  '

# Row 17
$ws.Range("D17").Value = '2.285.51'
$ws.Range("E17").Value = '  +0.34%  '

# Row 18
$ws.Range("D18").Value = '42.656.11'
$ws.Range("E18").Value = '  +1.37%  '

# Row 19
$ws.Range("D19").Value = '''7.26'
$ws.Range("E19").Value = '  -0.76%  '

# Row 20
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '''0.0000104'
$ws.Range("E20").Value = '  -0.25%  '

# Row 21
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").Value = '''13.52'
$ws.Range("E21").Value = '  +6.79%  '

# Row 22
$ws.Range("D22").Value = '''72.94'
$ws.Range("E22").Value = '  +0.53%  '

# Row 23
$ws.Range("E23").Value = '  -2.73%  '

# Row 24
$ws.Range("D24").Value = '''262.24'
$ws.Range("E24").Value = '  -1.66%  '

# Row 25
$ws.Range("D25").Value = '''2.17'
$ws.Range("E25").Value = '  -1.44%  '

# Row 26
$ws.Range("E26").Value = '  +0.36%  '

# Row 27
$ws.Range("D27").Value = '''10.67'
$ws.Range("E27").Value = '  -0.72%  '

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '''2.32'
$ws.Range("E28").Value = '  -0.02%  '

# Row 29
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = '''6.96'
$ws.Range("E29").Value = '  +15.80%  '

# Row 30
$ws.Range("D30").Value = '''22.23'
$ws.Range("E30").Value = '  -0.39%  '

# Row 31
$ws.Range("D31").Value = '''164.95'
$ws.Range("E31").Value = '  +0.43%  '

# Row 32
$ws.Range("D32").Value = '''35.51'
$ws.Range("E32").Value = '  -4.02%  '

# Row 33
$ws.Range("D33").Value = '''0.0857'
$ws.Range("E33").Value = '  -0.80%  '

# Row 34
$ws.Range("E34").Value = '  -1.71%  '

# Row 35
$ws.Range("E35").Value = '  +0.79%  '

# Row 36
$ws.Range("D36").Value = '''0.111'
$ws.Range("E36").Value = '  -2.33%  '

# Row 37
$ws.Range("D37").Value = '''4.48'
$ws.Range("E37").Value = '  -1.44%  '

# Row 38
$ws.Range("E38").Value = '  -1.01%  '

# Row 39
$ws.Range("D39").Value = '''3.71'
$ws.Range("E39").Value = '  +1.77%  '

# Row 40
$ws.Range("E40").Value = '  -0.95%  '

# Row 41
$ws.Range("D41").Value = '''1.57'
$ws.Range("E41").Value = '  +4.29%  '

# Row 42
$ws.Range("D42").Value = '''98.67'
$ws.Range("E42").Value = '  +7.91%  '

# Row 43
$ws.Range("D43").Value = '''68.82'
$ws.Range("E43").Value = '  +1.32%  '

# Row 44
$ws.Range("E44").Value = '  +0.48%  '

# Row 45
$ws.Range("D45").Value = '''0.225'
$ws.Range("E45").Value = '  +0.88%  '

# Row 46
$ws.Range("D46").Value = '1.721.87'
$ws.Range("E46").Value = '  +7.23%  '

# Row 47
$ws.Range("D47").Value = '''11.93'
$ws.Range("E47").Value = '  +0.65%  '

# Row 48
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '''110.08'
$ws.Range("E48").Value = '  -4.30%  '

# Row 49
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").Value = '''77.81'
$ws.Range("E49").Value = '  -1.08%  '

# Row 50
$ws.Range("D50").Value = '''5.16'
$ws.Range("E50").Value = '  -0.48%  '

# Row 51
$ws.Range("D51").Value = '''8.61'
$ws.Range("E51").Value = '  -3.33%  '
